$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(0.0668,0.0646,0.0601,0.0606,0.0603,0.0596,0.0588,0.0583,0.058,0.058,0.0584,0.0582,0.0583,0.0583,0.0586,0.0587,0.0593,0.0609,0.0609,0.062,0.062,0.0617,0.0609,0.0603,0.0588,0.0578,0.057,0.0562,0.0569,0.0579,0.0592,0.0608,0.0626,0.0645,0.0668,0.0694,0.0723,0.0752,0.0781,0.0808,0.0839,0.0871,0.0917,0.0912,0.0936,0.0962,0.099,0.1021,0.1052,0.1081,0.1111,0.1149,0.1199,0.1251,0.1307,0.1362,0.1514,0.1466,0.1532,0.1577,0.1633,0.1689)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value = $values[$i]
    $cell.NumberFormat = "0.00%"
}

$ws.Range("C2:C63").HorizontalAlignment = -4108
$ws.Range("C2:C63").VerticalAlignment = -4108

$ws.Range("C3").Select()
